# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking values must stay TEXT (matches source inlineStr cells)
# Force text number format before assigning so Excel does not auto-convert them to numbers.
$textCells = @('D5', 'D6', 'D8', 'D9', 'D15', 'D16', 'D19', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D29', 'D30', 'D34', 'D37', 'D38', 'D39', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D49', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '27.168.17'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('D3').Value = '1.651.26'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  -0.90%  '
$ws.Range('D5').Value = '219.83'
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('D6').Value = '0.502'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E7').Value = '  -0.77%  '
$ws.Range('D8').Value = '0.255'
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('D9').Value = '0.0627'
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('E10').Value = '  +2.36%  '
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').Value = '1.880.28'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').Value = '1.650.81'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').Value = '0.534'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').Value = '66.36'
$ws.Range('E16').Value = '  +2.37%  '
$ws.Range('D17').Value = '27.128.58'
$ws.Range('E17').Value = '  +1.09%  '
$ws.Range('D18').Value = '0.0₃0737'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = '223.27'
$ws.Range('E19').Value = '  +3.86%  '
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('D21').Value = '6.83'
$ws.Range('E21').Value = '  +8.89%  '
$ws.Range('D22').Value = '4.44'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').Value = '2.41'
$ws.Range('E23').Value = '  -2.71%  '
$ws.Range('D24').Value = '9.29'
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('D25').Value = '147.21'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').Value = '7.39'
$ws.Range('E27').Value = '  +2.57%  '
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').Value = '15.93'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').Value = '0.0514'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('E32').Value = '  +0.72%  '
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('D34').Value = '1.58'
$ws.Range('E34').Value = '  +2.54%  '
$ws.Range('D35').Value = '1.270.68'
$ws.Range('E35').Value = '  -1.96%  '
$ws.Range('E36').Value = '  -0.52%  '
$ws.Range('D37').Value = '0.0175'
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('D38').Value = '0.540'
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').Value = '0.828'
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('D41').Value = '0.805'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D42').Value = '5.39'
$ws.Range('E42').Value = '  +0.91%  '
$ws.Range('D43').Value = '1.790.74'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('D44').Value = '62.03'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').Value = '92.65'
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('D46').Value = '2.07'
$ws.Range('E46').Value = '  -7.69%  '
$ws.Range('D47').Value = '1.62'
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('E48').Value = '  -0.95%  '
$ws.Range('D49').Value = '7.66'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('D51').Value = '0.405'
$ws.Range('E51').Value = '  -0.54%  '
